$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.255.25"
$ws.Range("E2").Value = "  -3.33%  "
$ws.Range("D3").Value = "2.685.44"
$ws.Range("E3").Value = "  -7.55%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "498.38"
$ws.Range("E5").Value = "  -5.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.62"
$ws.Range("E6").Value = "  -3.09%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("E8").Value = "  -4.67%  "
$ws.Range("D9").Value = "2.695.90"
$ws.Range("E9").Value = "  -7.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.96"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.102"
$ws.Range("E11").Value = "  -5.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.340"
$ws.Range("E12").Value = "  -3.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.127"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "3.163.02"
$ws.Range("E14").Value = "  -7.32%  "
$ws.Range("D15").Value = "58.243.41"
$ws.Range("E15").Value = "  -3.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.12"
$ws.Range("E16").Value = "  -6.54%  "
$ws.Range("D17").Value = "2.691.95"
$ws.Range("E17").Value = "  -7.45%  "
$ws.Range("E18").Value = "  -5.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.62"
$ws.Range("E19").Value = "  -5.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.73"
$ws.Range("E20").Value = "  -6.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "336.59"
$ws.Range("E21").Value = "  -6.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.10"
$ws.Range("E22").Value = "  -7.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "61.74"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("E26").Value = "  -3.64%  "
$ws.Range("E27").Value = "  -7.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.994"
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").Value = "0.0₃0809"
$ws.Range("E29").Value = "  -5.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.28"
$ws.Range("E30").Value = "  -5.12%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  -4.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.75"
$ws.Range("E33").Value = "  -4.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "147.08"
$ws.Range("E34").Value = "  -3.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.08"
$ws.Range("E35").Value = "  -5.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.25"
$ws.Range("E36").Value = "  -5.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.910"
$ws.Range("E37").Value = "  -8.49%  "
$ws.Range("E38").Value = "  -7.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.83"
$ws.Range("E39").Value = "  -5.17%  "
$ws.Range("E40").Value = "  -6.42%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "2.141.70"
$ws.Range("E42").Value = "  -8.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.45"
$ws.Range("E43").Value = "  -5.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0547"
$ws.Range("E44").Value = "  -3.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.590"
$ws.Range("E45").Value = "  -8.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.34"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.39"
$ws.Range("E47").Value = "  -11.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.57"
$ws.Range("E48").Value = "  -5.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0222"
$ws.Range("E49").Value = "  -4.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0873"
$ws.Range("E50").Value = "  -5.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.39"
$ws.Range("E51").Value = "  -4.81%  "
